$wb = $excel.ActiveWorkbook

# The file "2ce35108-dbdb-4dfe-a285-e06b47185c04.md" received a new handoff,
# so its "Latest Handoff Date/Datetime" is refreshed on every sheet (row 5
# in each sheet corresponds to this file).

# Overview sheet: "Latest Handoff Date" column (D)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-03-24 20:49:51"

# zh-cn sheet: "Latest Handoff Datetime" column (E)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-24 20:49:44"

# de-de sheet: "Latest Handoff Datetime" column (E)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-24 20:49:51"
